$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped from this seed's sample ---
# Row 26 ("RM 232") is removed outright.
$ws.Rows.Item(26).Delete()
# Row 28 ("SC 92") was directly below it; after the first delete it is now
# row 27, so delete that next. Everything below shifts up accordingly,
# turning the former A1:F35 range into A1:F33.
$ws.Rows.Item(27).Delete()

# --- Cell-level value swaps (values imputed/removed for this seed) ---
$ws.Range("C3").Value = 11.2
$ws.Range("F4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("F12").Value = 17.45
$ws.Range("F15").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()

# --- Fill in values that are now missing for the rows that shifted up ---
# (former "SC 132" / "SC 193" rows, now rows 31 and 32 after the deletions)
$ws.Range("F31").Value = 17.18
$ws.Range("C32").Value = 10.5
$ws.Range("F32").Value = 17.39
